$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; A="all"; B="Intercept"; C="-43.8252333487047"; D="0.03285282928370578"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=3; A="all"; B="TrailerPublishYear"; C="0.02362618538139191"; D="0.02037410306100642"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=4; A="all"; B="TrailerPublishDays"; C="0.3891309111991351"; D="1.019440973731556e-23"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=5; A="all"; B="TrailerDuration"; C="0.006769402093940205"; D="4.288545333344304e-08"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=6; A="all"; B="ProductionBudget"; C="0.06353883292422918"; D="2.353804242694394e-21"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=7; A="all"; B="Action"; C="0.1218990873078701"; D="0.1223873254900664"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=8; A="all"; B="Comedy"; C="-0.02462753093527517"; D="0.7352825019571567"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=9; A="all"; B="Documentary"; C="-0.6714849375736641"; D="3.838930075153281e-12"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=10; A="all"; B="Drama"; C="-0.116107062659266"; D="0.09189021232048347"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=11; A="all"; B="PG-13"; C="0.4717506508983828"; D="9.619713999603922e-07"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=12; A="all"; B="R"; C="0.5312788849591763"; D="1.796318178140144e-09"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=13; A="all"; B="Not Rated"; C="-0.3002566725303573"; D="0.0008178513884850787"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=14; A="all"; B="TopStars"; C="0.4604889734022516"; D="3.621381104698225e-13"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=15; A="all"; B="TopStudios"; C="0.740535880491475"; D="1.643218334867517e-17"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=16; A="all"; B="ProductionCountry"; C="0.1716332612165648"; D="0.01177912543282537"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=17; A="all"; B="RatioFaceNo"; C="-1.012376435441653"; D="0.01425548237235368"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=18; A="all"; B="FaceNo"; C="0.07035145876140922"; D="0.008625544325999168"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=19; A="all"; B="AvgFaceSize"; C="-0.008812991728604653"; D="0.3211047679285141"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=20; A="all"; B="RatioFaceCoverage"; C="2.339502210442216"; D="0.5413028639566368"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=21; A="all"; B="RatioFemale"; C="-0.04788207345458586"; D="0.8280903833998621"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=22; A="all"; B="AverageAge"; C="-0.1069903726951685"; D="9.058533046428861e-12"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=23; A="all"; B="RatioSad"; C="0.06733083444275073"; D="0.8774929985221944"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=24; A="all"; B="RatioHappy"; C="-0.7121430658814929"; D="0.1846024448225281"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=25; A="all"; B="RatioFear"; C="0.6661142237900112"; D="0.2386291171570335"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=26; A="all"; B="RatioAngry"; C="0.7402135422236131"; D="0.1439308837834031"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=27; A="all"; B="RatioSurprise"; C="-0.110647385359129"; D="0.9774868025276635"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=28; A="all"; B="RatioDisgust"; C="230.4003120814403"; D="0.2514330054216104"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=29; A="all"; B="RatioAsian"; C="0.1483252544924223"; D="0.6331107975981328"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=30; A="all"; B="RatioIndian"; C="0.09220587816716463"; D="0.9892179180608173"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=31; A="all"; B="RatioBlack"; C="-1.064911857949007"; D="0.4277533955382481"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=32; A="all"; B="RatioMiddle"; C="2.487039094084162"; D="0.05988240896943365"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=33; A="all"; B="RatioHispanic"; C="0.8099952688434624"; D="0.579017715115594"; E="0.4412257720869083"; F="0.4337658233077069"; G="1.216707847767516"; H="0.8807848215939169"},
    @{Row=34; A="cont"; B="Intercept"; C="-40.98164797904721"; D="0.04162042136235278"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=35; A="cont"; B="TrailerPublishYear"; C="0.02030642471571845"; D="0.0416704482270633"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=36; A="cont"; B="TrailerPublishDays"; C="0.3900458181078471"; D="1.720484303159662e-23"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=37; A="cont"; B="TrailerDuration"; C="0.008232302159308843"; D="5.327214001419852e-16"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=38; A="cont"; B="ProductionBudget"; C="0.06859348883087343"; D="1.201541555368325e-24"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=39; A="cont"; B="Action"; C="0.1492002496941142"; D="0.0526721806325547"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=40; A="cont"; B="Comedy"; C="-0.02010186390894747"; D="0.766070962337832"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=41; A="cont"; B="Documentary"; C="-0.6582390124252219"; D="6.719795663874087e-13"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=42; A="cont"; B="Drama"; C="-0.1030218627816445"; D="0.1324395214722167"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=43; A="cont"; B="PG-13"; C="0.4946225307248153"; D="3.662215626930061e-07"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=44; A="cont"; B="R"; C="0.5400133408287348"; D="1.352137861142645e-09"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=45; A="cont"; B="Not Rated"; C="-0.3236089438349208"; D="0.000339126620250968"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=46; A="cont"; B="TopStars"; C="0.4188153081035941"; D="4.359940702255215e-11"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=47; A="cont"; B="TopStudios"; C="0.7557257309017842"; D="5.742208729394474e-18"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"},
    @{Row=48; A="cont"; B="ProductionCountry"; C="0.1484980106232207"; D="0.03009809895157389"; E="0.4226810072822857"; F="0.4192254853079171"; G="1.257088308686512"; H="0.8959622673844697"}
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Range($ws.Cells.Item($r, 3), $ws.Cells.Item($r, 8)).Style = "Normal"
}
